$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C holds text labels that look like numbers in scientific notation
# (e.g. "6.7e-01"). Force them to stay text so Excel doesn't coerce them
# into numeric cells, then restore the default ("Normal") cell style so no
# extra formatting is applied.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 3) "6.5e-01"
Set-TextValue $ws.Cells.Item(3, 3) "5.6e-01"
Set-TextValue $ws.Cells.Item(4, 3) "1.9e+01"

# Row 2 (ethyl benzene)
$ws.Range("B2").Value = 0.6544732734249999
$ws.Range("P2").Value = 180
$ws.Range("Q2").Value = 337
$ws.Range("V2").Value = 514485
$ws.Range("W2").Value = 3776824
$ws.Range("X2").Value = 34.1321878108
$ws.Range("Y2").Value = -86.8429036877

# Row 3 (toluene)
$ws.Range("B3").Value = 0.555684854795
$ws.Range("P3").Value = 180
$ws.Range("Q3").Value = 337
$ws.Range("V3").Value = 514485
$ws.Range("W3").Value = 3776824
$ws.Range("X3").Value = 34.1321878108
$ws.Range("Y3").Value = -86.8429036877

# Row 4 (xylenes (mixed))
$ws.Range("B4").Value = 19.2143474236
$ws.Range("P4").Value = 180
$ws.Range("Q4").Value = 337
$ws.Range("V4").Value = 514485
$ws.Range("W4").Value = 3776824
$ws.Range("X4").Value = 34.1321878108
$ws.Range("Y4").Value = -86.8429036877
